$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "9"
$ws.Range("B5").Value = "[BUG] deployment failing"
$ws.Range("C5").Value = "open"
$ws.Range("D5").Value = "2025-03-24T08:23:49Z"
$ws.Range("E5").Value = "bug"
